# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Leve profit-tracking sheets
# as captured by the upstream scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1014.2857
$ws.Range("H88").Value = 1274.2727
$ws.Range("J88").Value = 1188.8334
$ws.Range("L88").Value = 1188.8334
$ws.Range("N88").Value = -2000.8334
$ws.Range("H91").Value = 1274.2727
$ws.Range("J91").Value = 1188.8334
$ws.Range("L91").Value = 1188.8334
$ws.Range("N91").Value = -3996.8334
$ws.Range("H98").Value = 3842.5
$ws.Range("I98").Value = 3689.4666
$ws.Range("K98").Value = 3689.4666
$ws.Range("M98").Value = -2191.4666
$ws.Range("H113").Value = 14298.875
$ws.Range("J113").Value = 5748.5
$ws.Range("L113").Value = 5748.5
$ws.Range("N113").Value = -12256.5
$ws.Range("H122").Value = 3842.5
$ws.Range("I122").Value = 3689.4666
$ws.Range("K122").Value = 11068.3998
$ws.Range("M122").Value = -8618.399800000001
$ws.Range("H138").Value = 4148.9814
$ws.Range("I138").Value = 1537.4286
$ws.Range("J138").Value = 4537.936
$ws.Range("K138").Value = 4612.2858
$ws.Range("L138").Value = 13613.808
$ws.Range("M138").Value = 527.7142000000003
$ws.Range("N138").Value = -23893.808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1528.4667
$ws.Range("I45").Value = 1536.4445
$ws.Range("J45").Value = 1516.5
$ws.Range("K45").Value = 1536.4445
$ws.Range("L45").Value = 1516.5
$ws.Range("M45").Value = -1159.4445
$ws.Range("N45").Value = -2270.5
$ws.Range("H92").Value = 39999.5
$ws.Range("J92").Value = 39999.5
$ws.Range("L92").Value = 39999.5
$ws.Range("N92").Value = -44991.5
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 2475.5557
$ws.Range("I122").Value = 2313.5
$ws.Range("K122").Value = 6940.5
$ws.Range("M122").Value = -4490.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1217.3
$ws.Range("I22").Value = 1478
$ws.Range("K22").Value = 1478
$ws.Range("M22").Value = -1305
$ws.Range("H86").Value = 2116.3215
$ws.Range("I86").Value = 2211.85
$ws.Range("K86").Value = 2211.85
$ws.Range("M86").Value = -1088.85
$ws.Range("H89").Value = 2116.3215
$ws.Range("I89").Value = 2211.85
$ws.Range("K89").Value = 11059.25
$ws.Range("M89").Value = -5443.25
$ws.Range("H94").Value = 902.9474
$ws.Range("I94").Value = 1101.16
$ws.Range("J94").Value = 521.7692
$ws.Range("K94").Value = 1101.16
$ws.Range("L94").Value = 521.7692
$ws.Range("M94").Value = -650.1600000000001
$ws.Range("N94").Value = -1423.7692
$ws.Range("H134").Value = 61775.477
$ws.Range("I134").Value = 100777
$ws.Range("J134").Value = 26319.545
$ws.Range("K134").Value = 302331
$ws.Range("L134").Value = 78958.63499999999
$ws.Range("M134").Value = -299796
$ws.Range("N134").Value = -84028.63499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 517.5
$ws.Range("I25").Value = 517.5
$ws.Range("K25").Value = 517.5
$ws.Range("M25").Value = -343.5
$ws.Range("H29").Value = 5000
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 5000
$ws.Range("N29").Value = -5586
$ws.Range("H31").Value = 7911.25
$ws.Range("I31").Value = 1273.0952
$ws.Range("K31").Value = 1273.0952
$ws.Range("M31").Value = -978.0952
$ws.Range("H34").Value = 7911.25
$ws.Range("I34").Value = 1273.0952
$ws.Range("K34").Value = 1273.0952
$ws.Range("M34").Value = -1071.0952
$ws.Range("H122").Value = 2127.5
$ws.Range("I122").Value = 2127.5
$ws.Range("K122").Value = 6382.5
$ws.Range("M122").Value = -3932.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5500
$ws.Range("I62").Value = 5500
$ws.Range("K62").Value = 16500
$ws.Range("M62").Value = -15814
$ws.Range("H65").Value = 5500
$ws.Range("I65").Value = 5500
$ws.Range("K65").Value = 49500
$ws.Range("M65").Value = -46068
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H81").Value = 8500633
$ws.Range("I81").Value = 766.6667
$ws.Range("J81").Value = 17000500
$ws.Range("K81").Value = 2300.0001
$ws.Range("L81").Value = 51001500
$ws.Range("M81").Value = -1177.0001
$ws.Range("N81").Value = -51003746
$ws.Range("H84").Value = 8500633
$ws.Range("I84").Value = 766.6667
$ws.Range("J84").Value = 17000500
$ws.Range("K84").Value = 6900.0003
$ws.Range("L84").Value = 153004500
$ws.Range("M84").Value = -1284.0003
$ws.Range("N84").Value = -153015732
$ws.Range("H98").Value = 900
$ws.Range("J98").Value = 900
$ws.Range("L98").Value = 2700
$ws.Range("N98").Value = -5696

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 2012000
$ws.Range("H80").Value = 15345
$ws.Range("I80").Value = 15345
$ws.Range("K80").Value = 15345
$ws.Range("M80").Value = -14347
$ws.Range("H83").Value = 15345
$ws.Range("I83").Value = 15345
$ws.Range("K83").Value = 76725
$ws.Range("M83").Value = -71733
$ws.Range("H122").Value = 4200.9165
$ws.Range("I122").Value = 4436.8887
$ws.Range("K122").Value = 13310.6661
$ws.Range("M122").Value = -10860.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2730.2273
$ws.Range("I46").Value = 1057.2858
$ws.Range("K46").Value = 1057.2858
$ws.Range("M46").Value = -869.2858000000001
$ws.Range("H55").Value = 1851.5238
$ws.Range("I55").Value = 2289.8572
$ws.Range("J55").Value = 1632.3572
$ws.Range("K55").Value = 2289.8572
$ws.Range("L55").Value = 1632.3572
$ws.Range("M55").Value = -2116.8572
$ws.Range("N55").Value = -1978.3572
$ws.Range("H64").Value = 97324.5
$ws.Range("J64").Value = 97324.5
$ws.Range("L64").Value = 97324.5
$ws.Range("N64").Value = -97774.5
$ws.Range("H67").Value = 97324.5
$ws.Range("J67").Value = 97324.5
$ws.Range("L67").Value = 97324.5
$ws.Range("N67").Value = -98884.5
$ws.Range("H122").Value = 9865.706
$ws.Range("I122").Value = 11345.875
$ws.Range("K122").Value = 34037.625
$ws.Range("M122").Value = -31587.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 18333.334
$ws.Range("I31").Value = 15000
$ws.Range("J31").Value = 20000
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = -14652
$ws.Range("N31").Value = -20696
$ws.Range("H62").Value = 13174
$ws.Range("I62").Value = 10513.6
$ws.Range("J62").Value = 16499.5
$ws.Range("K62").Value = 10513.6
$ws.Range("L62").Value = 16499.5
$ws.Range("M62").Value = -9889.6
$ws.Range("N62").Value = -17747.5
$ws.Range("H63").Value = 30185.5
$ws.Range("J63").Value = 30185.5
$ws.Range("L63").Value = 30185.5
$ws.Range("N63").Value = -31433.5
$ws.Range("H65").Value = 13174
$ws.Range("I65").Value = 10513.6
$ws.Range("J65").Value = 16499.5
$ws.Range("K65").Value = 52568
$ws.Range("L65").Value = 82497.5
$ws.Range("M65").Value = -49448
$ws.Range("N65").Value = -88737.5
$ws.Range("H66").Value = 30185.5
$ws.Range("J66").Value = 30185.5
$ws.Range("L66").Value = 90556.5
$ws.Range("N66").Value = -96796.5
$ws.Range("H68").Value = 43385
$ws.Range("J68").Value = 43385
$ws.Range("L68").Value = 43385
$ws.Range("N68").Value = -45007
$ws.Range("H69").Value = 29904.2
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 29904.2
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 29904.2
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -31402.2
$ws.Range("H71").Value = 43385
$ws.Range("J71").Value = 43385
$ws.Range("L71").Value = 130155
$ws.Range("N71").Value = -138267
$ws.Range("H72").Value = 29904.2
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 29904.2
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 89712.60000000001
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -97200.60000000001
$ws.Range("H80").Value = 32651
$ws.Range("J80").Value = 32651
$ws.Range("L80").Value = 32651
$ws.Range("N80").Value = -34647
$ws.Range("H81").Value = 1396.8182
$ws.Range("I81").Value = 1408.1875
$ws.Range("J81").Value = 1366.5
$ws.Range("K81").Value = 2816.375
$ws.Range("L81").Value = 2733
$ws.Range("M81").Value = -1755.375
$ws.Range("N81").Value = -4855
$ws.Range("H83").Value = 32651
$ws.Range("J83").Value = 32651
$ws.Range("L83").Value = 97953
$ws.Range("N83").Value = -107937
$ws.Range("H84").Value = 1396.8182
$ws.Range("I84").Value = 1408.1875
$ws.Range("J84").Value = 1366.5
$ws.Range("K84").Value = 14081.875
$ws.Range("L84").Value = 13665
$ws.Range("M84").Value = -8777.875
$ws.Range("N84").Value = -24273

